$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 1.455362044514542
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 3.537761648806719
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 7.143138311642302

# Row 3
$ws.Range("B3").Value = 1.455362044514542
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 3.537761648806719
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 7.143138311642302

# Row 4
$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 3.537761648806719
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("G4").Value = 8.974608811992548
